$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header cells: "<label>_old" -> "<label>_FV2404", "<label>_new" -> "<label>_FV2410"
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $val = $cell.Value()
    $cell.Value = $val -replace "_old$", "_FV2404"
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $val = $cell.Value()
    $cell.Value = $val -replace "_new$", "_FV2410"
}

# 2. Turn the data range into a native Excel Table ("Table1") spanning A1:U72.
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U72"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# 3. Freeze the header row (split/freeze pane after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
